$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 2
$ws.Range("D2").Value = "'29.189.74"
$ws.Range("D2").Style = "Normal"
$ws.Range("E2").Value = "  -0.69%  "

# Row 3
$ws.Range("D3").Value = "'1.860.82"
$ws.Range("D3").Style = "Normal"
$ws.Range("E3").Value = "  -0.92%  "

# Row 4
$ws.Range("D4").Value = "'0.9996"
$ws.Range("D4").Style = "Normal"
$ws.Range("E4").Value = "  -0.15%  "

# Row 5
$ws.Range("B5").Value = "XRP"
$ws.Range("C5").Value = "https://coinranking.com/coin/-l8Mn2pVlRs-p+xrp-xrp"
$ws.Range("D5").Value = "'0.7024"
$ws.Range("D5").Style = "Normal"
$ws.Range("E5").Value = "  -1.91%  "

# Row 6
$ws.Range("B6").Value = "BNB"
$ws.Range("C6").Value = "https://coinranking.com/coin/WcwrkfNI4FUAe+bnb-bnb"
$ws.Range("D6").Value = "'241.75"
$ws.Range("D6").Style = "Normal"
$ws.Range("E6").Value = "  -0.80%  "

# Row 7
$ws.Range("D7").Value = "'0.9998"
$ws.Range("D7").Style = "Normal"

# Row 8
$ws.Range("D8").Value = "'0.07820"
$ws.Range("D8").Style = "Normal"
$ws.Range("E8").Value = "  -1.99%  "

# Row 9
$ws.Range("E9").Value = "  -1.18%  "

# Row 10
$ws.Range("D10").Value = "'23.88"
$ws.Range("D10").Style = "Normal"
$ws.Range("E10").Value = "  -4.08%  "

# Row 11
$ws.Range("D11").Value = "'0.07802"
$ws.Range("D11").Style = "Normal"
$ws.Range("E11").Value = "  -3.46%  "

# Row 12
$ws.Range("B12").Value = "Litecoin"
$ws.Range("C12").Value = "https://coinranking.com/coin/D7B1x_ks7WhV5+litecoin-ltc"
$ws.Range("D12").Value = "'92.60"
$ws.Range("D12").Style = "Normal"
$ws.Range("E12").Value = "  -2.19%  "

# Row 13
$ws.Range("B13").Value = "Polkadot"
$ws.Range("C13").Value = "https://coinranking.com/coin/25W7FG7om+polkadot-dot"
$ws.Range("D13").Value = "'5.122"
$ws.Range("D13").Style = "Normal"
$ws.Range("E13").Value = "  -1.99%  "

# Row 14
$ws.Range("B14").Value = "Polygon"
$ws.Range("C14").Value = "https://coinranking.com/coin/uW2tk-ILY0ii+polygon-matic"
$ws.Range("D14").Value = "'0.6913"
$ws.Range("D14").Style = "Normal"
$ws.Range("E14").Value = "  -2.29%  "

# Row 15
$ws.Range("B15").Value = "WrappedEther"
$ws.Range("C15").Value = "https://coinranking.com/coin/Mtfb0obXVh59u+wrappedether-weth"
$ws.Range("D15").Value = "'1.788.22"
$ws.Range("D15").Style = "Normal"
$ws.Range("E15").Value = "  -4.62%  "

# Row 16
$ws.Range("D16").Value = "'6.552"
$ws.Range("D16").Style = "Normal"
$ws.Range("E16").Value = "  +2.46%  "

# Row 17
$ws.Range("D17").Value = "'0.000008447"
$ws.Range("D17").Style = "Normal"
$ws.Range("E17").Value = "  +0.05%  "

# Row 18
$ws.Range("D18").Value = "'29.206.06"
$ws.Range("D18").Style = "Normal"
$ws.Range("E18").Value = "  -0.65%  "

# Row 19
$ws.Range("D19").Value = "'250.24"
$ws.Range("D19").Style = "Normal"
$ws.Range("E19").Value = "  -1.03%  "

# Row 20
$ws.Range("D20").Value = "'2.111.18"
$ws.Range("D20").Style = "Normal"
$ws.Range("E20").Value = "  -1.06%  "

# Row 21
$ws.Range("E21").Value = "  -3.19%  "

# Row 22
$ws.Range("D22").Value = "'1.0000"
$ws.Range("D22").Style = "Normal"
$ws.Range("E22").Value = "  -0.06%  "

# Row 23
$ws.Range("D23").Value = "'7.591"
$ws.Range("D23").Style = "Normal"
$ws.Range("E23").Value = "  -1.03%  "

# Row 24
$ws.Range("E24").Value = "  -0.13%  "

# Row 25
$ws.Range("D25").Value = "'0.1532"
$ws.Range("D25").Style = "Normal"
$ws.Range("E25").Value = "  -2.77%  "

# Row 26
$ws.Range("D26").Value = "'160.47"
$ws.Range("D26").Style = "Normal"
$ws.Range("E26").Value = "  -0.91%  "

# Row 27
$ws.Range("D27").Value = "'8.887"
$ws.Range("D27").Style = "Normal"
$ws.Range("E27").Value = "  -1.97%  "

# Row 28
$ws.Range("E28").Value = "  -2.29%  "

# Row 29
$ws.Range("D29").Value = "'1.570"
$ws.Range("D29").Style = "Normal"
$ws.Range("E29").Value = "  +4.08%  "

# Row 30
$ws.Range("E30").Value = "  -3.20%  "

# Row 31
$ws.Range("E31").Value = "  -1.51%  "

# Row 32
$ws.Range("D32").Value = "'1.213"
$ws.Range("D32").Style = "Normal"
$ws.Range("E32").Value = "  -0.99%  "

# Row 33
$ws.Range("D33").Value = "'0.05222"
$ws.Range("D33").Style = "Normal"
$ws.Range("E33").Value = "  -1.63%  "

# Row 34
$ws.Range("B34").Value = "LidoDAOToken"
$ws.Range("C34").Value = "https://coinranking.com/coin/Pe93bIOD2+lidodaotoken-ldo"
$ws.Range("D34").Value = "'1.875"
$ws.Range("D34").Style = "Normal"
$ws.Range("E34").Value = "  -3.48%  "

# Row 35
$ws.Range("B35").Value = "ImmutableX"
$ws.Range("C35").Value = "https://coinranking.com/coin/Z96jIvLU7+immutablex-imx"
$ws.Range("D35").Value = "'0.7572"
$ws.Range("D35").Style = "Normal"
$ws.Range("E35").Value = "  -0.02%  "

# Row 36
$ws.Range("E36").Value = "  +0.15%  "

# Row 37
$ws.Range("D37").Value = "'2.710"
$ws.Range("D37").Style = "Normal"
$ws.Range("E37").Value = "  +0.31%  "

# Row 38
$ws.Range("E38").Value = "  -1.28%  "

# Row 39
$ws.Range("D39").Value = "'1.223.19"
$ws.Range("D39").Style = "Normal"
$ws.Range("E39").Value = "  -4.07%  "

# Row 40
$ws.Range("E40").Value = "  -1.24%  "

# Row 41
$ws.Range("D41").Value = "'0.9001"
$ws.Range("D41").Style = "Normal"

# Row 42
$ws.Range("D42").Value = "'110.21"
$ws.Range("D42").Style = "Normal"
$ws.Range("E42").Value = "  -1.21%  "

# Row 43
$ws.Range("D43").Value = "'5.816"
$ws.Range("D43").Style = "Normal"
$ws.Range("E43").Value = "  -9.30%  "

# Row 44
$ws.Range("D44").Value = "'0.9993"
$ws.Range("D44").Style = "Normal"
$ws.Range("E44").Value = "  -0.15%  "

# Row 45
$ws.Range("D45").Value = "'2.008.80"
$ws.Range("D45").Style = "Normal"
$ws.Range("E45").Value = "  -1.04%  "

# Row 46
$ws.Range("E46").Value = "  -4.61%  "

# Row 47
$ws.Range("D47").Value = "'64.96"
$ws.Range("D47").Style = "Normal"
$ws.Range("E47").Value = "  -12.35%  "

# Row 48
$ws.Range("D48").Value = "'0.5186"
$ws.Range("D48").Style = "Normal"
$ws.Range("E48").Value = "  -0.42%  "

# Row 49
$ws.Range("D49").Value = "'9.511"
$ws.Range("D49").Style = "Normal"
$ws.Range("E49").Value = "  -0.11%  "

# Row 50
$ws.Range("D50").Value = "'1.770"
$ws.Range("D50").Style = "Normal"
$ws.Range("E50").Value = "  -2.00%  "

# Row 51
$ws.Range("D51").Value = "'7.033"
$ws.Range("D51").Style = "Normal"
$ws.Range("E51").Value = "  -0.89%  "
